$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bugs")

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Alta"
$ws.Range("C4").Value = "Exclusão de Entidade de Facturação"
$ws.Range("D4").Value = "Em aberto"

$ws.Range("A4:B4").HorizontalAlignment = -4108

$ws.Range("D5").Select()
